$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "01‏/05‏/2025 02:08:51 م"
$ws.Range("B3").Value = "NRC"
$ws.Range("C3").Value = "C4"
$ws.Range("D3").Value = "الرحلة 2"
$ws.Range("E3").Value = "ايتا"
$ws.Range("F3").Value = "احمد"

# "222" looks numeric - force text storage (quote-prefix), then drop the
# resulting "quote prefix" style so the cell keeps the default style.
$ws.Range("G3").Value = "'222"
$ws.Range("G3").Style = "Normal"

# Empty text cell (matches the other blank-but-text cells in this sheet,
# e.g. H2) - same quote-prefix trick with nothing after the apostrophe.
$ws.Range("H3").Value = "'"
$ws.Range("H3").Style = "Normal"
